$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 05:12"

# --- Row 6: Estados Unidos (refreshed case numbers) ---
$ws.Range("B6").Value = 68421
$ws.Range("C6").Value = 210
$ws.Range("E6").Value = 66995
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 1032

# --- Row 20: Australia (refreshed case numbers) ---
$ws.Range("B20").Value = 2728
$ws.Range("C20").Value = 52
$ws.Range("E20").Value = 2598
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 12

# --- Rows 44/45: India overtakes Rusia in total cases, so the two
#     countries swap places in the (descending, sorted-by-total-cases)
#     table. India's figures are also refreshed with newer numbers.
$ws.Range("A44").Value = "India"
$ws.Range("B44").Value = 664
$ws.Range("C44").Value = 7
$ws.Range("D44").Value = 43
$ws.Range("E44").Value = 609
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 12

$ws.Range("A45").Value = "Rusia"
$ws.Range("B45").Value = 658
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 29
$ws.Range("E45").Value = 626
$ws.Range("F45").Value = 8
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 3
